$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update shared string text used in B43
$ws.Range("B43").Value = "Максимальная сумма, руб."

# 2. Update H25:H38 formulas (simplify IF(SUM(...)) to IF(Gn<=Fn, 0, Gn-Fn))
for ($r = 25; $r -le 38; $r++) {
    $ws.Range("H$r").Formula = "=IF(G$r<=F$r, 0, G$r-F$r)"
}

# 3. Update D35:D38 formulas
$ws.Range("D35").Formula = "=A1*1.1/2"
$ws.Range("D36").Formula = "=D35"
$ws.Range("D37").Formula = "=D36"
$ws.Range("D38").Formula = "=D37"

# 4. Update C43 formula
$ws.Range("C43").Formula = "=MAX(K3:K38)"
